# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before "总计", populated with
#    the fund-holding breakdown for the quarter (mirrors the layout of the
#    other quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名).
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet and renumber the
#    existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" sheet, inserted just before "总计"
# ---------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")

# Clone the most recent fully-populated quarter sheet so the new sheet
# inherits the same sheet-level formatting (outline props, margins,
# header style, border style, column layout) instead of Excel's blank
# worksheet defaults.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template has 6 data rows; 2022-Q1 only has 5, so drop the extra one.
$newSheet.Rows.Item(7).Delete()

$fundRows = @(
    @("010671", "景顺长城大中华混合(QDII)美元", "10.35", "82.59", "8.20", "0.8487", 1),
    @("262001", "景顺长城大中华混合(QDII)", "10.35", "82.59", "8.20", "0.8487", 1),
    @("000927", "博时大中华亚太精选股票(QDII) - 美元现汇", "0.32", "92.94", "5.03", "0.0161", 6),
    @("050015", "博时大中华亚太精选股票(QDII) -人民币", "0.32", "92.94", "5.03", "0.0161", 6),
    @("040021", "华安大中华升级股票(QDII)", "0.26", "87.37", "5.06", "0.0132", 1)
)

# Columns B:G hold text (fund codes with leading zeros, and numeric-looking
# figures that must stay strings) - force text format before writing so the
# COM layer doesn't coerce them to numbers, then drop back to the sheet's
# normal style so no stray per-cell formatting is left behind.
$textRange = $newSheet.Range("B2:G6")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

$textRange.Style = "Normal"

# ---------------------------------------------------------------------
# Part 2: prepend a "2022-Q1" row to "总计" and renumber the index column
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# Pick up the correct formatting for the new row (index-column style,
# plain data-column style) from the row just below it.
$ws.Range("A3:D3").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "2022-Q1"
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = 1.74

# Renumber the 0-based index column for the rows that shifted down.
for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Restore the workbook's original active sheet (the sheet copy and edits
# above shift focus onto the new/changed sheets).
$wb.Worksheets.Item("2020-Q4").Activate()
